# Upload / add the third ("Random") uniform scenario sheet, rename the
# existing two sheets, refresh the "Constant" sheet's last data row and
# make the new sheet the active one.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the two existing sheets -------------------------------
$wsConstant     = $wb.Worksheets.Item(1)
$wsProportional = $wb.Worksheets.Item(2)

$wsConstant.Name     = "Uniform - Constant"
$wsProportional.Name = "Uniform - Proportional"

# --- 2. Fix up the last data row on the "Proportional" sheet ----------
# (row 11 now holds real simulation output instead of placeholder zeros)
$wsProportional.Range("D11").Value = 36852
$wsProportional.Range("E11").Value = 147306

# Selection on that sheet is no longer the "active" one - point it at the
# full used range instead of a single cell.
$wsProportional.Range("A1:E12").Select() | Out-Null

# --- 3. Add the new "Random" sheet ------------------------------------
# Copy the "Proportional" sheet (keeps formulas/column widths/formats)
# to the end of the workbook, then overwrite its data with the random
# scenario numbers.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsProportional.Copy([Type]::Missing, $lastSheet) | Out-Null

$wsRandom = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsRandom.Name = "Uniform - Random"

$randomData = @(
    @(44651, 29494, 75829),
    @(44688,     0,     0),
    @(    0,     0,     0),
    @(    0,     0,     0),
    @(    0,     0,     0),
    @(    0,     0,     0),
    @(    0,     0,     0),
    @(    0,     0,     0),
    @(    0,     0,     0),
    @(    0,     0,     0)
)

for ($i = 0; $i -lt $randomData.Count; $i++) {
    $row = 2 + $i
    $vals = $randomData[$i]
    $wsRandom.Range("C$row").Value = $vals[0]
    $wsRandom.Range("D$row").Value = $vals[1]
    $wsRandom.Range("E$row").Value = $vals[2]
}

# D12/E12 already carry =AVERAGE(...) formulas copied from the source
# sheet - they will recalc automatically.

$wsRandom.Range("D3").Select() | Out-Null

# --- 4. Make the new "Random" sheet the active tab --------------------
$wsRandom.Activate()

$wb.Application.Calculate()
